# This workbook tracks daily price observations for "Piña" (pineapple) at
# the Terminal Hortofrutícola Agro Chillán market. The commit adds one new
# daily observation. The new record is inserted at row 115 (pushing the
# existing rows 115-217 down by one, to 116-218); its "shape" columns
# (Calidad/Volumen/Precio mínimo/Precio máximo/Precio promedio ponderado/
# Unidad de comercialización/Precio $/Kg/Kg por unidad) are copied from the
# row immediately above it (what was row 114), while the date (Fecha)
# receives a new value, consistent with the rest of the workbook where a
# newly-logged row reuses the preceding row's commercial figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 115; everything at/after 115 shifts down one.
$ws.Rows("115:115").Insert()

# Populate the newly inserted row 115 with its data.
$ws.Range("A115").Value = 7
$ws.Range("B115").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C115").Value = "Ñuble"
$ws.Range("D115").Value = 44669
$ws.Range("E115").Value = 16
$ws.Range("F115").Value = "Fruta"
$ws.Range("G115").Value = 100108
$ws.Range("H115").Value = "Tropicales y subtropicales"
$ws.Range("I115").Value = 100108005
$ws.Range("J115").Value = "Piña"
$ws.Range("K115").Value = "Caramelo"
$ws.Range("L115").Value = "Segunda"
$ws.Range("M115").Value = 60
$ws.Range("N115").Value = 16000
$ws.Range("O115").Value = 17000
$ws.Range("P115").Value = 16500
$ws.Range("Q115").Value = "$/caja 14 unidades"
$ws.Range("R115").Value = "Ecuador"
$ws.Range("S115").Value = 1179
$ws.Range("T115").Value = 14
